$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.151145577430725
$ws.Range("B1").Value = 1.326534748077393
$ws.Range("C1").Value = 1.672534227371216
$ws.Range("D1").Value = 3.377491235733032
$ws.Range("E1").Value = 15
